$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.310.38"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.881.32"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.35%  "
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.41"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.355"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.46"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0971"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.29%  "
$ws.Range("D14").Value = "2.157.58"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.765"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.90"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "1.886.32"
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("D18").Value = "35.385.44"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").Value = "0.0₃0819"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "243.60"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.77"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("E23").Value = "  -1.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.65"
$ws.Range("D24").ClearFormats()
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.27"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.54"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -1.83%  "
$ws.Range("D31").Value = "4.128.45"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +9.60%  "
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.26"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  -3.03%  "
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  -3.26%  "
$ws.Range("E40").Value = "  +10.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.45"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0216"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.98"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -6.98%  "
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "1.303.75"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("E46").Value = "  -2.49%  "
$ws.Range("E47").Value = "  +6.92%  "
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.72"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "12.11"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.21"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.57%  "
